$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 313, shifting rows 313:328 down to 314:329
$ws.Rows.Item(313).Insert()

# Populate the new row 313 with the new data record
$ws.Cells.Item(313, 1).Value = 4
$ws.Cells.Item(313, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(313, 3).Value = "Los Lagos"
$ws.Cells.Item(313, 4).Value = 44753
$ws.Cells.Item(313, 5).Value = 10
$ws.Cells.Item(313, 6).Value = 100112045
$ws.Cells.Item(313, 7).Value = "Zapallo"
$ws.Cells.Item(313, 8).Value = "Paine"
$ws.Cells.Item(313, 9).Value = "1a (guarda)"
$ws.Cells.Item(313, 10).Value = 500
$ws.Cells.Item(313, 11).Value = 500
$ws.Cells.Item(313, 12).Value = 500
$ws.Cells.Item(313, 13).Value = 500
$ws.Cells.Item(313, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(313, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(313, 16).Value = 500
$ws.Cells.Item(313, 17).Value = 1
$ws.Cells.Item(313, 18).Value = "Hortaliza"
